# Copy and paste (1 to 1) the Descrimination Percent statistics for each
# country from HIVDescrimination.xls into column C of GlobalHIVMerged.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

${ws}.Range('C2').Value = '60.2  (Source: DHS 2015 )'
${ws}.Range('C3').Value = '64.1  (Source: 2008 DHS)'
${ws}.Range('C4').Value = '51  (Source: INCAPSIDA 2010)'
${ws}.Range('C5').Value = '73  (Source: 2010 DHS)'
${ws}.Range('C6').Value = '30.8  (Source: Knowledge Attitudes and Practices Study 2014,)'
${ws}.Range('C7').Value = '55.4  (Source: DHS 2011-2012)'
${ws}.Range('C8').Value = '44.7  (Source: 2008-09 DHS)'
${ws}.Range('C9').Value = '57.1  (Source: 2011-12 MICS)'
${ws}.Range('C10').Value = '13.2  (Source: BAIS IV)'
${ws}.Range('C11').Value = '62.3  (Source: 2010 DHS)'
${ws}.Range('C12').Value = '25.5  (Source: 2010 DHS)'
${ws}.Range('C13').Value = '19.4  (Source: DHS 2014)'
${ws}.Range('C14').Value = '40.5  (Source: 2011 DHS)'
${ws}.Range('C15').Value = '32.6  (Source: 2010 MICS)'
${ws}.Range('C16').Value = '55.8  (Source: 2012 DHS)'
${ws}.Range('C17').Value = '35.4  (Source: 2011-12 DHS)'
${ws}.Range('C18').Value = '16.4  (Source: 2014 MICS)'
${ws}.Range('C19').Value = '44.5  (Source: 2011-12 DHS)'
${ws}.Range('C20').Value = '49.2  (Source: 2013-14 DHS )'
${ws}.Range('C21').Value = '49.3  (Source: 2013 DHS)'
${ws}.Range('C22').Value = '59.9  (Source: 2011 DHS)'
${ws}.Range('C23').Value = '25.3  (Source: 2012 DHS)'
${ws}.Range('C24').Value = '51  (Source: 2013 DHS)'
${ws}.Range('C25').Value = '67.7  (Source: 2014 DHS)'
${ws}.Range('C26').Value = '60.6  (Source: Encuesta Nacional de Salud Materno Infantil 2008-09)'
${ws}.Range('C27').Value = '80.1  (Source: 2012 DHS)'
${ws}.Range('C28').Value = '57.7  (Source: DHS 2012)'
${ws}.Range('C29').Value = '44.9  (Source: 2011-12 DHS)'
${ws}.Range('C30').Value = '62.8  (Source: 2012 DHS)'
${ws}.Range('C31').Value = '71  (Source: Knowledge Attitudes, Behaviors and Practices 2012,)'
${ws}.Range('C32').Value = '64.8  (Source: 2010-11 MICS)'
${ws}.Range('C33').Value = '11.9  (Source: 2014 DHS)'
${ws}.Range('C34').Value = '57.2  (Source: 2012 DHS)'
${ws}.Range('C35').Value = '53.5  (Source: 2011-12 MICS)'
${ws}.Range('C36').Value = '13.9  (Source: DHS 2014)'
${ws}.Range('C37').Value = '52.7  (Source: 2013 DHS)'
${ws}.Range('C38').Value = '14.9  (Source: DHS 2015-2016)'
${ws}.Range('C39').Value = '45.8  (Source: 2012-13 DHS)'
${ws}.Range('C40').Value = '75.6  (Source: 2010 MICS)'
${ws}.Range('C41').Value = '52  (Source: 2013 MICS)'
${ws}.Range('C42').Value = '28  (Source: 2011 DHS)'
${ws}.Range('C43').Value = '13  (Source: DHS 2013)'
${ws}.Range('C44').Value = '28.3  (Source: 2011 DHS)'
${ws}.Range('C45').Value = '71.5  (Source: DHS 2012)'
${ws}.Range('C46').Value = '46.8  (Source: 2013 DHS)'
${ws}.Range('C47').Value = '49  (Source: 2012-13 DHS)'
${ws}.Range('C48').Value = '70.8  (Source: 2012 MICS)'
${ws}.Range('C49').Value = '9.9  (Source: DHS 2014-2015)'
${ws}.Range('C50').Value = '84.2  (Source: OECS BSS 2005-06)'
${ws}.Range('C51').Value = '51.7  (Source: 2014 DHS)'
${ws}.Range('C52').Value = '53.4  (Source: 2013 DHS)'
${ws}.Range('C53').Value = '12.5  (Source: 2010 MICS)'
${ws}.Range('C54').Value = '63.2  (Source: 2009-10 DHS)'
${ws}.Range('C55').Value = '45.8  (Source: 2013-14 DHS )'
${ws}.Range('C56').Value = '26.2  (Source: 2011 DHS)'
${ws}.Range('C57').Value = '65.1  (Source: 2012 MICS)'
${ws}.Range('C58').Value = '26.8  (Source: 2010 DHS)'
${ws}.Range('C59').Value = '18  (Source: 2013-14 DHS )'
${ws}.Range('C60').Value = '17.8  (Source: DHS 2015 )'

# Update the saved selection / window position to match the authored state.
$ws.Range("H10").Select()

$win = $wb.Windows.Item(1)
$win.Top = 2820
